$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1081
$ws.Range("I18").Value = 1079.1
$ws.Range("J18").Value = 1100
$ws.Range("K18").Value = 1079.1
$ws.Range("L18").Value = 1100
$ws.Range("M18").Value = -795.0999999999999
$ws.Range("N18").Value = -1668
$ws.Range("H111").Value = 536.25
$ws.Range("I111").Value = 535.5
$ws.Range("J111").Value = 537.375
$ws.Range("K111").Value = 1606.5
$ws.Range("L111").Value = 1612.125
$ws.Range("M111").Value = 1460.5
$ws.Range("N111").Value = -7746.125
$ws.Range("H137").Value = 29413624
$ws.Range("I137").Value = 1186.069
$ws.Range("J137").Value = 200005760
$ws.Range("K137").Value = 3558.207
$ws.Range("L137").Value = 600017280
$ws.Range("M137").Value = -1008.207
$ws.Range("N137").Value = -600022380

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 9573.429
$ws.Range("J37").Value = 11402.8
$ws.Range("L37").Value = 11402.8
$ws.Range("N37").Value = -11948.8
$ws.Range("H55").Value = 18350.75
$ws.Range("J55").Value = 18350.75
$ws.Range("L55").Value = 18350.75
$ws.Range("N55").Value = -18980.75
$ws.Range("H68").Value = 52184.75
$ws.Range("J68").Value = 54765
$ws.Range("L68").Value = 54765
$ws.Range("N68").Value = -56387
$ws.Range("H71").Value = 52184.75
$ws.Range("J71").Value = 54765
$ws.Range("L71").Value = 164295
$ws.Range("N71").Value = -172407
$ws.Range("H80").Value = 20391.25
$ws.Range("J80").Value = 20391.25
$ws.Range("L80").Value = 20391.25
$ws.Range("N80").Value = -22387.25
$ws.Range("H83").Value = 20391.25
$ws.Range("J83").Value = 20391.25
$ws.Range("L83").Value = 61173.75
$ws.Range("N83").Value = -71157.75

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 808.25806
$ws.Range("I20").Value = 726.94116
$ws.Range("J20").Value = 907
$ws.Range("K20").Value = 726.94116
$ws.Range("L20").Value = 907
$ws.Range("M20").Value = -479.94116
$ws.Range("N20").Value = -1401
$ws.Range("H80").Value = 263.6
$ws.Range("I80").Value = 54
$ws.Range("J80").Value = 316
$ws.Range("K80").Value = 54
$ws.Range("L80").Value = 316
$ws.Range("M80").Value = 944
$ws.Range("N80").Value = -2312
$ws.Range("H82").Value = 18255.092
$ws.Range("I82").Value = 10439.25
$ws.Range("J82").Value = 22721.285
$ws.Range("K82").Value = 10439.25
$ws.Range("L82").Value = 22721.285
$ws.Range("M82").Value = -10056.25
$ws.Range("N82").Value = -23487.285
$ws.Range("H83").Value = 263.6
$ws.Range("I83").Value = 54
$ws.Range("J83").Value = 316
$ws.Range("K83").Value = 270
$ws.Range("L83").Value = 1580
$ws.Range("M83").Value = 4722
$ws.Range("N83").Value = -11564
$ws.Range("H85").Value = 18255.092
$ws.Range("I85").Value = 10439.25
$ws.Range("J85").Value = 22721.285
$ws.Range("K85").Value = 10439.25
$ws.Range("L85").Value = 22721.285
$ws.Range("M85").Value = -9113.25
$ws.Range("N85").Value = -25373.285
$ws.Range("H107").Value = 914.9167
$ws.Range("I107").Value = 805.3125
$ws.Range("J107").Value = 1134.125
$ws.Range("K107").Value = 805.3125
$ws.Range("L107").Value = 1134.125
$ws.Range("M107").Value = 1114.6875
$ws.Range("N107").Value = -4974.125
$ws.Range("H140").Value = 34852.668
$ws.Range("J140").Value = 34852.668
$ws.Range("L140").Value = 34852.668
$ws.Range("N140").Value = -45212.668

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 400.25
$ws.Range("I22").Value = 200.125
$ws.Range("J22").Value = 800.5
$ws.Range("K22").Value = 200.125
$ws.Range("L22").Value = 800.5
$ws.Range("M22").Value = 149.875
$ws.Range("N22").Value = -1500.5
$ws.Range("H31").Value = 1624
$ws.Range("I31").Value = 1300.3334
$ws.Range("K31").Value = 1300.3334
$ws.Range("M31").Value = -1005.3334
$ws.Range("H34").Value = 1624
$ws.Range("I34").Value = 1300.3334
$ws.Range("K34").Value = 1300.3334
$ws.Range("M34").Value = -1098.3334
$ws.Range("H50").Value = 10197.333
$ws.Range("J50").Value = 10539.429
$ws.Range("L50").Value = 10539.429
$ws.Range("N50").Value = -11789.429
$ws.Range("H51").Value = 10754.637
$ws.Range("J51").Value = 10920.1
$ws.Range("L51").Value = 10920.1
$ws.Range("N51").Value = -12392.1
$ws.Range("H58").Value = 2292.2886
$ws.Range("I58").Value = 1029.826
$ws.Range("J58").Value = 3293.5518
$ws.Range("K58").Value = 1029.826
$ws.Range("L58").Value = 3293.5518
$ws.Range("M58").Value = -826.826
$ws.Range("N58").Value = -3699.5518
$ws.Range("H60").Value = 8525.75
$ws.Range("J60").Value = 10741.2
$ws.Range("L60").Value = 10741.2
$ws.Range("N60").Value = -11763.2
$ws.Range("H61").Value = 10754.637
$ws.Range("J61").Value = 10920.1
$ws.Range("L61").Value = 10920.1
$ws.Range("N61").Value = -11616.1
$ws.Range("H68").Value = 19348.75
$ws.Range("J68").Value = 19348.75
$ws.Range("L68").Value = 19348.75
$ws.Range("N68").Value = -20846.75
$ws.Range("H71").Value = 19348.75
$ws.Range("J71").Value = 19348.75
$ws.Range("L71").Value = 58046.25
$ws.Range("N71").Value = -65534.25
$ws.Range("H74").Value = 16567.875
$ws.Range("J74").Value = 16567.875
$ws.Range("L74").Value = 16567.875
$ws.Range("N74").Value = -18315.875
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H77").Value = 16567.875
$ws.Range("J77").Value = 16567.875
$ws.Range("L77").Value = 49703.625
$ws.Range("N77").Value = -58439.625
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H99").Value = 1579.75
$ws.Range("I99").Value = 1306
$ws.Range("J99").Value = 1853.5
$ws.Range("K99").Value = 1306
$ws.Range("L99").Value = 1853.5
$ws.Range("M99").Value = 192
$ws.Range("N99").Value = -4849.5
$ws.Range("H126").Value = 1579.75
$ws.Range("I126").Value = 1306
$ws.Range("J126").Value = 1853.5
$ws.Range("K126").Value = 3918
$ws.Range("L126").Value = 5560.5
$ws.Range("M126").Value = -1448
$ws.Range("N126").Value = -10500.5
$ws.Range("H129").Value = 49999
$ws.Range("J129").Value = 49999
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999
$ws.Range("H132").Value = 1992.7106
$ws.Range("I132").Value = 1897.0333
$ws.Range("K132").Value = 5691.0999
$ws.Range("M132").Value = -3161.0999
$ws.Range("H136").Value = 2292.2886
$ws.Range("I136").Value = 1029.826
$ws.Range("J136").Value = 3293.5518
$ws.Range("K136").Value = 3089.478
$ws.Range("L136").Value = 9880.6554
$ws.Range("M136").Value = -539.4780000000001
$ws.Range("N136").Value = -14980.6554
$ws.Range("H141").Value = 31642.117
$ws.Range("J141").Value = 31642.117
$ws.Range("L141").Value = 31642.117
$ws.Range("N141").Value = -42002.117

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 58826430
$ws.Range("I34").Value = 333.75
$ws.Range("J34").Value = 76926776
$ws.Range("K34").Value = 1001.25
$ws.Range("L34").Value = 230780328
$ws.Range("M34").Value = -917.25
$ws.Range("N34").Value = -230780496
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H55").Value = 98.388885
$ws.Range("I55").Value = 96.77778000000001
$ws.Range("J55").Value = 100
$ws.Range("K55").Value = 290.33334
$ws.Range("L55").Value = 300
$ws.Range("M55").Value = -113.33334
$ws.Range("N55").Value = -654

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 8118
$ws.Range("J15").Value = 8118
$ws.Range("L15").Value = 8118
$ws.Range("N15").Value = -8694
$ws.Range("H80").Value = 76340.664
$ws.Range("I80").Value = 3263.125
$ws.Range("J80").Value = 159857.86
$ws.Range("K80").Value = 3263.125
$ws.Range("L80").Value = 159857.86
$ws.Range("M80").Value = -2265.125
$ws.Range("N80").Value = -161853.86
$ws.Range("H81").Value = 8118
$ws.Range("J81").Value = 8118
$ws.Range("L81").Value = 8118
$ws.Range("N81").Value = -10114
$ws.Range("H83").Value = 76340.664
$ws.Range("I83").Value = 3263.125
$ws.Range("J83").Value = 159857.86
$ws.Range("K83").Value = 16315.625
$ws.Range("L83").Value = 799289.2999999999
$ws.Range("M83").Value = -11323.625
$ws.Range("N83").Value = -809273.2999999999
$ws.Range("H84").Value = 8118
$ws.Range("J84").Value = 8118
$ws.Range("L84").Value = 24354
$ws.Range("N84").Value = -34338
$ws.Range("H122").Value = 3081.9546
$ws.Range("I122").Value = 1798.8276
$ws.Range("J122").Value = 5562.6665
$ws.Range("K122").Value = 5396.4828
$ws.Range("L122").Value = 16687.9995
$ws.Range("M122").Value = -2946.4828
$ws.Range("N122").Value = -21587.9995

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3125.5
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 5251
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 5251
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -5627
$ws.Range("H82").Value = 1357
$ws.Range("I82").Value = 1052.4445
$ws.Range("J82").Value = 1813.8334
$ws.Range("K82").Value = 1052.4445
$ws.Range("L82").Value = 1813.8334
$ws.Range("M82").Value = -691.4445000000001
$ws.Range("N82").Value = -2535.8334
$ws.Range("H85").Value = 1357
$ws.Range("I85").Value = 1052.4445
$ws.Range("J85").Value = 1813.8334
$ws.Range("K85").Value = 1052.4445
$ws.Range("L85").Value = 1813.8334
$ws.Range("M85").Value = 195.5554999999999
$ws.Range("N85").Value = -4309.8334

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()
$ws.Range("H140").Value = 22538.166
$ws.Range("J140").Value = 22538.166
$ws.Range("L140").Value = 22538.166
$ws.Range("N140").Value = -32898.166
